$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update startDate value (D2): 2024-09-01 -> 2025-01-01
$ws.Range("D2").Value = "2025-01-01"

# Update currency value (E2): USD -> AED
$ws.Range("E2").Value = "AED"

# Update conversion value (C2): 3.6725 -> 1
$ws.Range("C2").Value = 1

# Update active cell selection to C2
$ws.Range("C2").Select()
